$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date value in A1 (45735 -> 45749)
$ws.Range("A1").Value = 45749

# Move the active selection from D11 to A3
$ws.Range("A3").Select()
